# Scheduled-runner update: refresh cached marketboard-derived profit figures
# (currentAveragePrice* / Leve*Price* / Leve*Profit* columns) across every class
# table in the Odin_Profits workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3494.1667
$ws.Range("I9").Value = 3494.1667
$ws.Range("K9").Value = 3494.1667
$ws.Range("M9").Value = -3325.1667
$ws.Range("H33").Value = 524.6923
$ws.Range("I33").Value = 485.08334
$ws.Range("K33").Value = 485.08334
$ws.Range("M33").Value = -256.08334
$ws.Range("H40").Value = 3431.3333
$ws.Range("I40").Value = 2300
$ws.Range("K40").Value = 2300
$ws.Range("M40").Value = -2125
$ws.Range("H125").Value = 2028.9231
$ws.Range("J125").Value = 2309.375
$ws.Range("L125").Value = 20784.375
$ws.Range("N125").Value = -25704.375
$ws.Range("H141").Value = 2875
$ws.Range("I141").Value = 2166.6667
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 6500.000100000001
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -1320.000100000001
$ws.Range("N141").Value = -25360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1986688.6
$ws.Range("I32").Value = 2106.843
$ws.Range("J32").Value = 22229422
$ws.Range("K32").Value = 2106.843
$ws.Range("L32").Value = 22229422
$ws.Range("M32").Value = -1819.843
$ws.Range("N32").Value = -22229996
$ws.Range("H61").Value = 5794.357
$ws.Range("I61").Value = 50000
$ws.Range("J61").Value = 4716.171
$ws.Range("K61").Value = 50000
$ws.Range("L61").Value = 4716.171
$ws.Range("M61").Value = -49788
$ws.Range("N61").Value = -5140.171
$ws.Range("H102").Value = 2254.2666
$ws.Range("I102").Value = 2190.0356
$ws.Range("K102").Value = 2190.0356
$ws.Range("M102").Value = -568.0356000000002
$ws.Range("H110").Value = 4791.147
$ws.Range("I110").Value = 2158.923
$ws.Range("J110").Value = 6420.619
$ws.Range("K110").Value = 2158.923
$ws.Range("L110").Value = 6420.619
$ws.Range("M110").Value = -113.9229999999998
$ws.Range("N110").Value = -10510.619
$ws.Range("H122").Value = 3287.8333
$ws.Range("I122").Value = 2602.5386
$ws.Range("K122").Value = 7807.6158
$ws.Range("M122").Value = -5357.6158
$ws.Range("H126").Value = 9980
$ws.Range("I126").Value = 9980
$ws.Range("K126").Value = 29940
$ws.Range("M126").Value = -27470
$ws.Range("H136").Value = 5794.357
$ws.Range("I136").Value = 50000
$ws.Range("J136").Value = 4716.171
$ws.Range("K136").Value = 150000
$ws.Range("L136").Value = 14148.513
$ws.Range("M136").Value = -147450
$ws.Range("N136").Value = -19248.513

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1115476.2
$ws.Range("I134").Value = 1319189
$ws.Range("J134").Value = 9607.286
$ws.Range("K134").Value = 3957567
$ws.Range("L134").Value = 28821.858
$ws.Range("M134").Value = -3955032
$ws.Range("N134").Value = -33891.858

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 36520.668
$ws.Range("J42").Value = 36520.668
$ws.Range("L42").Value = 36520.668
$ws.Range("N42").Value = -37706.668
$ws.Range("H58").Value = 41684200
$ws.Range("I58").Value = 83343610
$ws.Range("K58").Value = 83343610
$ws.Range("M58").Value = -83343407
$ws.Range("H92").Value = 22560.4
$ws.Range("J92").Value = 22560.4
$ws.Range("L92").Value = 22560.4
$ws.Range("N92").Value = -27552.4
$ws.Range("H136").Value = 41684200
$ws.Range("I136").Value = 83343610
$ws.Range("K136").Value = 250030830
$ws.Range("M136").Value = -250028280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1096.0526
$ws.Range("I5").Value = 705
$ws.Range("J5").Value = 1766.4286
$ws.Range("K5").Value = 2115
$ws.Range("L5").Value = 5299.2858
$ws.Range("M5").Value = -2003
$ws.Range("N5").Value = -5523.2858
$ws.Range("H49").Value = 4180
$ws.Range("I49").Value = 1966.6666
$ws.Range("K49").Value = 5899.9998
$ws.Range("M49").Value = -5743.9998
$ws.Range("H57").Value = 12599
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 12599
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 37797
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -38915
$ws.Range("H64").Value = 15456
$ws.Range("J64").Value = 15697.5
$ws.Range("L64").Value = 47092.5
$ws.Range("N64").Value = -47632.5
$ws.Range("H67").Value = 15456
$ws.Range("J67").Value = 15697.5
$ws.Range("L67").Value = 47092.5
$ws.Range("N67").Value = -48964.5
$ws.Range("H101").Value = 32503.666
$ws.Range("J101").Value = 35999.2
$ws.Range("L101").Value = 107997.6
$ws.Range("N101").Value = -112865.6
$ws.Range("H113").Value = 742.2941
$ws.Range("I113").Value = 417
$ws.Range("J113").Value = 877.8333
$ws.Range("K113").Value = 1251
$ws.Range("L113").Value = 2633.4999
$ws.Range("M113").Value = 919
$ws.Range("N113").Value = -6973.4999
$ws.Range("H122").Value = 5243.909
$ws.Range("I122").Value = 1116.75
$ws.Range("J122").Value = 6161.0557
$ws.Range("K122").Value = 10050.75
$ws.Range("L122").Value = 55449.5013
$ws.Range("M122").Value = -7600.75
$ws.Range("N122").Value = -60349.5013
$ws.Range("H135").Value = 1096.0526
$ws.Range("I135").Value = 705
$ws.Range("J135").Value = 1766.4286
$ws.Range("K135").Value = 6345
$ws.Range("L135").Value = 15897.8574
$ws.Range("M135").Value = -3810
$ws.Range("N135").Value = -20967.8574
$ws.Range("H136").Value = 100004800
$ws.Range("I136").Value = 83337336
$ws.Range("K136").Value = 250012008
$ws.Range("M136").Value = -250006908
$ws.Range("H137").Value = 3260.4119
$ws.Range("I137").Value = 1796.8
$ws.Range("J137").Value = 3870.25
$ws.Range("K137").Value = 5390.4
$ws.Range("L137").Value = 11610.75
$ws.Range("M137").Value = -290.3999999999996
$ws.Range("N137").Value = -21810.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10391.2
$ws.Range("I70").Value = 10489.5
$ws.Range("J70").Value = 9998
$ws.Range("K70").Value = 10489.5
$ws.Range("L70").Value = 9998
$ws.Range("M70").Value = -10219.5
$ws.Range("N70").Value = -10538
$ws.Range("H73").Value = 10391.2
$ws.Range("I73").Value = 10489.5
$ws.Range("J73").Value = 9998
$ws.Range("K73").Value = 10489.5
$ws.Range("L73").Value = 9998
$ws.Range("M73").Value = -9553.5
$ws.Range("N73").Value = -11870
$ws.Range("H80").Value = 4216.4136
$ws.Range("I80").Value = 3853.2727
$ws.Range("J80").Value = 4438.3335
$ws.Range("K80").Value = 3853.2727
$ws.Range("L80").Value = 4438.3335
$ws.Range("M80").Value = -2855.2727
$ws.Range("N80").Value = -6434.3335
$ws.Range("H83").Value = 4216.4136
$ws.Range("I83").Value = 3853.2727
$ws.Range("J83").Value = 4438.3335
$ws.Range("K83").Value = 19266.3635
$ws.Range("L83").Value = 22191.6675
$ws.Range("M83").Value = -14274.3635
$ws.Range("N83").Value = -32175.6675
$ws.Range("H113").Value = 9988.714
$ws.Range("I113").Value = 6555
$ws.Range("K113").Value = 6555
$ws.Range("M113").Value = -4385
$ws.Range("H132").Value = 32261930
$ws.Range("I132").Value = 55559524
$ws.Range("J132").Value = 3723.7693
$ws.Range("K132").Value = 166678572
$ws.Range("L132").Value = 11171.3079
$ws.Range("M132").Value = -166676042
$ws.Range("N132").Value = -16231.3079

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1919.8334
$ws.Range("I16").Value = 839.2222
$ws.Range("J16").Value = 5161.6665
$ws.Range("K16").Value = 839.2222
$ws.Range("L16").Value = 5161.6665
$ws.Range("M16").Value = -669.2222
$ws.Range("N16").Value = -5501.6665
$ws.Range("H55").Value = 4280
$ws.Range("J55").Value = 4985.769
$ws.Range("L55").Value = 4985.769
$ws.Range("N55").Value = -5331.769
$ws.Range("H61").Value = 7313.84
$ws.Range("I61").Value = 6407.533
$ws.Range("J61").Value = 8673.299999999999
$ws.Range("K61").Value = 6407.533
$ws.Range("L61").Value = 8673.299999999999
$ws.Range("M61").Value = -6205.533
$ws.Range("N61").Value = -9077.299999999999
$ws.Range("H113").Value = 7313.84
$ws.Range("I113").Value = 6407.533
$ws.Range("J113").Value = 8673.299999999999
$ws.Range("K113").Value = 6407.533
$ws.Range("L113").Value = 8673.299999999999
$ws.Range("M113").Value = -4237.533
$ws.Range("N113").Value = -13013.3
$ws.Range("H122").Value = 5174.357
$ws.Range("I122").Value = 4247.8
$ws.Range("J122").Value = 7490.75
$ws.Range("K122").Value = 12743.4
$ws.Range("L122").Value = 22472.25
$ws.Range("M122").Value = -10293.4
$ws.Range("N122").Value = -27372.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 21750
$ws.Range("J62").Value = 19550
$ws.Range("L62").Value = 19550
$ws.Range("N62").Value = -20798
$ws.Range("H65").Value = 21750
$ws.Range("J65").Value = 19550
$ws.Range("M65").Value = -103990
$ws.Range("H126").Value = 4815.222
$ws.Range("I126").Value = 2494.6667
$ws.Range("K126").Value = 7484.000100000001
$ws.Range("M126").Value = -5014.000100000001
